$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.915632
$ws.Range("H2").Value = 35.746896
$ws.Range("I2").Value = 0.2203762099850903
$ws.Range("J2").Value = 0.2203762099850904
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.09934133333334
$ws.Range("N2").Value = 63.29802400000001
$ws.Range("O2").Value = 0.2917236204149438
$ws.Range("P2").Value = 0.2917236204149438
$ws.Range("Q2").Value = 251.4119867703894
$ws.Range("R2").Value = 2262.707880933504
$ws.Range("S2").Value = 0.06428894583017443
$ws.Range("T2").Value = 0.06428894583017446

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.915632
$ws.Range("H3").Value = 35.746896
$ws.Range("I3").Value = 0.2203762099850903
$ws.Range("J3").Value = 0.2203762099850904
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.81943766666667
$ws.Range("N3").Value = 107.458313
$ws.Range("O3").Value = 0.4952465516465762
$ws.Range("P3").Value = 0.4952465516465762
$ws.Range("Q3").Value = 426.8112376829387
$ws.Range("R3").Value = 3841.301139146448
$ws.Range("S3").Value = 0.1091405580600577
$ws.Range("T3").Value = 0.1091405580600578

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.915632
$ws.Range("H4").Value = 35.746896
$ws.Range("I4").Value = 0.2203762099850903
$ws.Range("J4").Value = 0.2203762099850904
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 15.40769666666667
$ws.Range("N4").Value = 46.22309
$ws.Range("O4").Value = 0.2130298279384801
$ws.Range("P4").Value = 0.2130298279384801
$ws.Range("Q4").Value = 183.5924434476267
$ws.Range("R4").Value = 1652.33199102864
$ws.Range("S4").Value = 0.04694670609485815
$ws.Range("T4").Value = 0.04694670609485816

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.800487
$ws.Range("H5").Value = 89.40146100000001
$ws.Range("I5").Value = 0.5511514941691683
$ws.Range("J5").Value = 0.5511514941691684
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.09934133333334
$ws.Range("N5").Value = 63.29802400000001
$ws.Range("O5").Value = 0.2917236204149438
$ws.Range("P5").Value = 0.2917236204149438
$ws.Range("Q5").Value = 628.7706471125628
$ws.Range("R5").Value = 5658.935824013065
$ws.Range("S5").Value = 0.1607839092761355
$ws.Range("T5").Value = 0.1607839092761356

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.800487
$ws.Range("H6").Value = 89.40146100000001
$ws.Range("I6").Value = 0.5511514941691683
$ws.Range("J6").Value = 0.5511514941691684
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.81943766666667
$ws.Range("N6").Value = 107.458313
$ws.Range("O6").Value = 0.4952465516465762
$ws.Range("P6").Value = 0.4952465516465762
$ws.Range("Q6").Value = 1067.43668653281
$ws.Range("R6").Value = 9606.930178795295
$ws.Range("S6").Value = 0.2729558769221386
$ws.Range("T6").Value = 0.2729558769221387

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.800487
$ws.Range("H7").Value = 89.40146100000001
$ws.Range("I7").Value = 0.5511514941691683
$ws.Range("J7").Value = 0.5511514941691684
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 15.40769666666667
$ws.Range("N7").Value = 46.22309
$ws.Range("O7").Value = 0.2130298279384801
$ws.Range("P7").Value = 0.2130298279384801
$ws.Range("Q7").Value = 459.1568642149434
$ws.Range("R7").Value = 4132.411777934491
$ws.Range("S7").Value = 0.1174117079708941
$ws.Range("T7").Value = 0.1174117079708942

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.35338333333333
$ws.Range("H8").Value = 37.06015
$ws.Range("I8").Value = 0.2284722958457413
$ws.Range("J8").Value = 0.2284722958457413
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.09934133333334
$ws.Range("N8").Value = 63.29802400000001
$ws.Range("O8").Value = 0.2917236204149438
$ws.Range("P8").Value = 0.2917236204149438
$ws.Range("Q8").Value = 260.6482515715111
$ws.Range("R8").Value = 2345.8342641436
$ws.Range("S8").Value = 0.06665076530863376
$ws.Range("T8").Value = 0.06665076530863379

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.35338333333333
$ws.Range("H9").Value = 37.06015
$ws.Range("I9").Value = 0.2284722958457413
$ws.Range("J9").Value = 0.2284722958457413
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 35.81943766666667
$ws.Range("N9").Value = 107.458313
$ws.Range("O9").Value = 0.4952465516465762
$ws.Range("P9").Value = 0.4952465516465762
$ws.Range("Q9").Value = 442.4912442807722
$ws.Range("R9").Value = 3982.42119852695
$ws.Range("S9").Value = 0.1131501166643797
$ws.Range("T9").Value = 0.1131501166643798

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.35338333333333
$ws.Range("H10").Value = 37.06015
$ws.Range("I10").Value = 0.2284722958457413
$ws.Range("J10").Value = 0.2284722958457413
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.40769666666667
$ws.Range("N10").Value = 46.22309
$ws.Range("O10").Value = 0.2130298279384801
$ws.Range("P10").Value = 0.2130298279384801
$ws.Range("Q10").Value = 190.3371832070555
$ws.Range("R10").Value = 1713.0346488635
$ws.Range("S10").Value = 0.04867141387272778
$ws.Range("T10").Value = 0.04867141387272779

